$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label cell
$ws.Range("E9").Value = "x offset"

# New data block E10:G18 (row 14 intentionally left blank, matches source diff)
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0.001
$ws.Range("G10").Value = 0

$ws.Range("E11").Value = 0.001
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0

$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0.001

$ws.Range("E13").Value = -0.0137
$ws.Range("F13").Value = -0.0031
$ws.Range("G13").Value = 0.0021

$ws.Range("E15").Value = -0.000023453426853
$ws.Range("F15").Value = 0.001005282682668
$ws.Range("G15").Value = 0.000006025290218

$ws.Range("E16").Value = 0.001023621407377
$ws.Range("F16").Value = 0.000002512089703
$ws.Range("G16").Value = 0.000023397795619

$ws.Range("E17").Value = 0.000000329136254
$ws.Range("F17").Value = -0.00000381232572
$ws.Range("G17").Value = 0.00097939474317

$ws.Range("E18").Value = -0.01366720769154
$ws.Range("F18").Value = -0.003069754115766
$ws.Range("G18").Value = 0.002137111113079

# Update the selected cell to reflect the new active cell in the workbook view
$ws.Range("G23").Select()
